# Daily refresh of the "剩余" (days-remaining) tracker.
#
# Column D = 总天 (total cycle length, days)
# Column E = 剩余 (days remaining in the current cycle)
# Column F = 开始时间 (cycle start date, stored as an 8-digit yyyyMMdd number)
#
# Each day E counts down by 1. When a cycle's remaining-days counter would
# drop to 0, the cycle rolls over: E resets to the full cycle length (D)
# and F advances by D calendar days (the date the next cycle begins).
# Rows whose F value isn't a well-formed 8-digit yyyyMMdd date are left
# untouched (stale/bad data that the refresh routine skips).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)
    $fCell = $ws.Cells.Item($r, 6)

    $dVal = $dCell.Value2
    $eVal = $eCell.Value2
    $fVal = $fCell.Value2

    if ($dVal -eq $null -or $eVal -eq $null -or $fVal -eq $null) {
        continue
    }

    $fText = [string]([int]$fVal)
    if ($fText.Length -ne 8) {
        # Malformed start date (e.g. "202510929") - skip, don't touch.
        continue
    }

    $year = [int]$fText.Substring(0, 4)
    $month = [int]$fText.Substring(4, 2)
    $day = [int]$fText.Substring(6, 2)

    $parsedOk = $true
    try {
        $startDate = Get-Date -Year $year -Month $month -Day $day
    } catch {
        $parsedOk = $false
    }
    if (-not $parsedOk) {
        continue
    }

    $totalDays = [int]$dVal
    $remaining = [int]$eVal

    if ($remaining -le 1) {
        # Cycle rolls over: reset remaining to the full duration and push
        # the start date forward by that many days.
        $newRemaining = $totalDays
        $newStart = $startDate.AddDays($totalDays)
        $eCell.Value2 = $newRemaining
        $fCell.Value2 = [int]$newStart.ToString("yyyyMMdd")
    } else {
        $eCell.Value2 = $remaining - 1
    }
}
